# Release: Splash screen for 0.3
#
# 1) Every "Date Placeholder" shape's cached datetimeFigureOut field
#    text moves from 7/7/2010 -> 8/3/2010 (slide master + all 11
#    slide layouts).
# 2) The "Release 0.2.0" textbox on slide 1 becomes "Release 0.3.0".

$p = $ppt.ActivePresentation

function Update-DatePlaceholder {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "7/7/2010") {
                $tr.Text = "8/3/2010"
            }
        }
    }
}

# Slide master.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout off the master.
$layouts = $master.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    $layout = $layouts.Item($L)
    Update-DatePlaceholder $layout.Shapes
}

# Slide 1: "Release 0.2.0" -> "Release 0.3.0" (only the version run).
$s = $p.Slides.Item(1)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $tr = $shp.TextFrame.TextRange
        $full = $tr.Text
        $idx = $full.IndexOf("0.2.0")
        if ($idx -ge 0) {
            $sub = $tr.Characters($idx + 1, 5)
            $sub.Text = "0.3.0"
        }
    }
}
